# AutoDocGenR.xlsx: remove the "RuleSet = Sequential / TRUE" row from the
# Investor Approval rule-table sheet. Deleting the entire row 2 shifts every
# row below it up by one (old row 11 -> new row 10), which is exactly the
# content/layout seen in the target workbook - no other cell content changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(2).Delete()

# Match the saved selection/cursor position left in the edited workbook.
$ws.Range("B13").Select()
